$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in "NA" for the duplicate_image_filename column (E) for data rows 2 through 21
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
